# Generate Report for Handoff
# A new handoff was generated for file
# "5685cf38-4c69-4098-a2ff-8993427d9e74" in both the zh-cn and de-de
# localization tables. Update the "Latest Handoff Datetime" column (H)
# for that file's row (row 4) in each language worksheet.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-10-14 07:31:07"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-10-14 07:31:18"
